# Auto-generated Excel COM-interop script
# Applies the "MAJ automatique BRVM via GitHub Actions" data refresh:
#  - Sheet "Recommandations": row 50 is removed (data now spans A1:G49),
#    and all data rows (2-49) are refreshed with the latest recommandation figures
#    (values recomputed, several titles reordered by the automatic ranking).
#  - Sheet "Top_YTD": the Progression YTD (%) values (column B) are refreshed
#    for the existing 10 ranked rows (A2:B11).

$wb = $excel.ActiveWorkbook
$wsReco = $wb.Worksheets.Item("Recommandations")
$wsYtd  = $wb.Worksheets.Item("Top_YTD")

# ---------------------------------------------------------------------------
# 1) "Recommandations" sheet: drop the last row (50) so the sheet shrinks back
#    to 49 rows (A1:G49), then rewrite rows 2-49 with the refreshed dataset.
# ---------------------------------------------------------------------------
$wsReco.Rows.Item(50).Delete()

$data1 = New-Object 'object[,]' 48,7
$data1[0,0] = "SUCRIVOIRE"
$data1[0,1] = 0
$data1[0,2] = 4
$data1[0,3] = 3920
$data1[0,4] = 970
$data1[0,5] = "🟡 Observer"
$data1[0,6] = "➖ Neutre"
$data1[1,0] = "BRVM - SERVICES PUBLICS"
$data1[1,1] = 0
$data1[1,2] = 8
$data1[1,3] = 3348.29
$data1[1,4] = 110.43
$data1[1,5] = "🟡 Observer"
$data1[1,6] = "➖ Neutre"
$data1[2,0] = "SAFCA CI"
$data1[2,1] = 0
$data1[2,2] = 4
$data1[2,3] = 2770
$data1[2,4] = 700
$data1[2,5] = "🟡 Observer"
$data1[2,6] = "➖ Neutre"
$data1[3,0] = "CFAO MOTORS CI"
$data1[3,1] = 0
$data1[3,2] = 4
$data1[3,3] = 2700
$data1[3,4] = 675
$data1[3,5] = "🟡 Observer"
$data1[3,6] = "➖ Neutre"
$data1[4,0] = "BRVM - AUTRES SECTEURS"
$data1[4,1] = 0
$data1[4,2] = 4
$data1[4,3] = 2656.63
$data1[4,4] = 659.84
$data1[4,5] = "🟡 Observer"
$data1[4,6] = "➖ Neutre"
$data1[5,0] = "NEI-CEDA CI"
$data1[5,1] = 0
$data1[5,2] = 4
$data1[5,3] = 2375
$data1[5,4] = 595
$data1[5,5] = "🟡 Observer"
$data1[5,6] = "➖ Neutre"
$data1[6,0] = "UNIWAX CI"
$data1[6,1] = 0
$data1[6,2] = 4
$data1[6,3] = 2340
$data1[6,4] = 580
$data1[6,5] = "🟡 Observer"
$data1[6,6] = "➖ Neutre"
$data1[7,0] = "SETAO CI"
$data1[7,1] = 0
$data1[7,2] = 4
$data1[7,3] = 2265
$data1[7,4] = 540
$data1[7,5] = "🟡 Observer"
$data1[7,6] = "➖ Neutre"
$data1[8,0] = "AIR LIQUIDE CI"
$data1[8,1] = 0
$data1[8,2] = 4
$data1[8,3] = 2135
$data1[8,4] = 530
$data1[8,5] = "🟡 Observer"
$data1[8,6] = "➖ Neutre"
$data1[9,0] = "BRVM - DISTRIBUTION"
$data1[9,1] = 0
$data1[9,2] = 4
$data1[9,3] = 1487.1
$data1[9,4] = 369.45
$data1[9,5] = "🟡 Observer"
$data1[9,6] = "➖ Neutre"
$data1[10,0] = "BRVM - TRANSPORT"
$data1[10,1] = 0
$data1[10,2] = 4
$data1[10,3] = 1401.33
$data1[10,4] = 348.8
$data1[10,5] = "🟡 Observer"
$data1[10,6] = "➖ Neutre"
$data1[11,0] = "BRVM - AGRICULTURE"
$data1[11,1] = 0
$data1[11,2] = 4
$data1[11,3] = 1304.75
$data1[11,4] = 315.62
$data1[11,5] = "🟡 Observer"
$data1[11,6] = "➖ Neutre"
$data1[12,0] = "BRVM - INDUSTRIE"
$data1[12,1] = 0
$data1[12,2] = 4
$data1[12,3] = 806.53
$data1[12,4] = 198.77
$data1[12,5] = "🟡 Observer"
$data1[12,6] = "➖ Neutre"
$data1[13,0] = "BRVM-PRINCIPAL"
$data1[13,1] = 0
$data1[13,2] = 4
$data1[13,3] = 705.0599999999999
$data1[13,4] = 174.45
$data1[13,5] = "🟡 Observer"
$data1[13,6] = "➖ Neutre"
$data1[14,0] = "BRVM - CONSOMMATION DE BASE"
$data1[14,1] = 0
$data1[14,2] = 4
$data1[14,3] = 704.48
$data1[14,4] = 172.36
$data1[14,5] = "🟡 Observer"
$data1[14,6] = "➖ Neutre"
$data1[15,0] = "BRVM - INDUSTRIELS"
$data1[15,1] = 0
$data1[15,2] = 4
$data1[15,3] = 533.67
$data1[15,4] = 132.48
$data1[15,5] = "🟡 Observer"
$data1[15,6] = "➖ Neutre"
$data1[16,0] = "BRVM-PRESTIGE"
$data1[16,1] = 0
$data1[16,2] = 4
$data1[16,3] = 521.75
$data1[16,4] = 129.37
$data1[16,5] = "🟡 Observer"
$data1[16,6] = "➖ Neutre"
$data1[17,0] = "BRVM - FINANCES"
$data1[17,1] = 0
$data1[17,2] = 4
$data1[17,3] = 486
$data1[17,4] = 120.74
$data1[17,5] = "🟡 Observer"
$data1[17,6] = "➖ Neutre"
$data1[18,0] = "BRVM - SERVICES FINANCIERS"
$data1[18,1] = 0
$data1[18,2] = 4
$data1[18,3] = 477.64
$data1[18,4] = 118.66
$data1[18,5] = "🟡 Observer"
$data1[18,6] = "➖ Neutre"
$data1[19,0] = "BRVM - ENERGIE"
$data1[19,1] = 0
$data1[19,2] = 4
$data1[19,3] = 442.01
$data1[19,4] = 109.47
$data1[19,5] = "🟡 Observer"
$data1[19,6] = "➖ Neutre"
$data1[20,0] = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$data1[20,1] = 0
$data1[20,2] = 4
$data1[20,3] = 424.51
$data1[20,4] = 106.23
$data1[20,5] = "🟡 Observer"
$data1[20,6] = "➖ Neutre"
$data1[21,0] = "BRVM - TELECOMMUNICATIONS"
$data1[21,1] = 0
$data1[21,2] = 4
$data1[21,3] = 377.51
$data1[21,4] = 93.66
$data1[21,5] = "🟡 Observer"
$data1[21,6] = "➖ Neutre"
$data1[22,0] = "BANK OF AFRICA ML (BOAM)"
$data1[22,1] = 2
$data1[22,2] = 0
$data1[22,3] = 11.78
$data1[22,4] = 4.99
$data1[22,5] = "🟡 Observer"
$data1[22,6] = "➖ Neutre"
$data1[23,0] = "TRACTAFRIC MOTORS CI (PRSC)"
$data1[23,1] = 2
$data1[23,2] = 1
$data1[23,3] = 6.15
$data1[23,4] = 6
$data1[23,5] = "🟡 Observer"
$data1[23,6] = "👀 À surveiller"
$data1[24,0] = "BANK OF AFRICA NG (BOAN)"
$data1[24,1] = 1
$data1[24,2] = 0
$data1[24,3] = 5.83
$data1[24,4] = 5.83
$data1[24,5] = "🟡 Observer"
$data1[24,6] = "➖ Neutre"
$data1[25,0] = "TOTALENERGIES MARKETING SN (TTLS)"
$data1[25,1] = 2
$data1[25,2] = 0
$data1[25,3] = 3.91
$data1[25,4] = 3.17
$data1[25,5] = "🟡 Observer"
$data1[25,6] = "➖ Neutre"
$data1[26,0] = "BERNABE CI (BNBC)"
$data1[26,1] = 2
$data1[26,2] = 2
$data1[26,3] = 3.55
$data1[26,4] = -2.69
$data1[26,5] = "🟡 Observer"
$data1[26,6] = "👀 À surveiller"
$data1[27,0] = "NSIA BANQUE COTE D'IVOIRE (NSBC)"
$data1[27,1] = 1
$data1[27,2] = 0
$data1[27,3] = 3.3
$data1[27,4] = 3.3
$data1[27,5] = "🟡 Observer"
$data1[27,6] = "➖ Neutre"
$data1[28,0] = "SAFCA CI (SAFC)"
$data1[28,1] = 1
$data1[28,2] = 0
$data1[28,3] = 3.08
$data1[28,4] = 3.08
$data1[28,5] = "🟡 Observer"
$data1[28,6] = "➖ Neutre"
$data1[29,0] = "VIVO ENERGY CI (SHEC)"
$data1[29,1] = 1
$data1[29,2] = 1
$data1[29,3] = 1.52
$data1[29,4] = 3.4
$data1[29,5] = "🟡 Observer"
$data1[29,6] = "👀 À surveiller"
$data1[30,0] = "CORIS BANK INTERNATIONAL (CBIBF)"
$data1[30,1] = 1
$data1[30,2] = 1
$data1[30,3] = 1.49
$data1[30,4] = -5.69
$data1[30,5] = "🟡 Observer"
$data1[30,6] = "👀 À surveiller"
$data1[31,0] = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$data1[31,1] = 1
$data1[31,2] = 0
$data1[31,3] = 1.4
$data1[31,4] = 1.4
$data1[31,5] = "🟡 Observer"
$data1[31,6] = "➖ Neutre"
$data1[32,0] = "FILTISAC CI (FTSC)"
$data1[32,1] = 1
$data1[32,2] = 0
$data1[32,3] = 0.88
$data1[32,4] = 0.88
$data1[32,5] = "🟡 Observer"
$data1[32,6] = "➖ Neutre"
$data1[33,0] = "UNIWAX CI (UNXC)"
$data1[33,1] = 1
$data1[33,2] = 1
$data1[33,3] = 0.51
$data1[33,4] = 7.41
$data1[33,5] = "🟡 Observer"
$data1[33,6] = "👀 À surveiller"
$data1[34,0] = "SICABLE CI (CABC)"
$data1[34,1] = 1
$data1[34,2] = 1
$data1[34,3] = 0.27
$data1[34,4] = 5.31
$data1[34,5] = "🟡 Observer"
$data1[34,6] = "👀 À surveiller"
$data1[35,0] = "TOTAL"
$data1[35,1] = 0
$data1[35,2] = 4
$data1[35,3] = 0
$data1[35,4] = 0
$data1[35,5] = "🟡 Observer"
$data1[35,6] = "➖ Neutre"
$data1[36,0] = "NEI-CEDA CI (NEIC)"
$data1[36,1] = 1
$data1[36,2] = 1
$data1[36,3] = -0.74
$data1[36,4] = 2.59
$data1[36,5] = "🟡 Observer"
$data1[36,6] = "👀 À surveiller"
$data1[37,0] = "ECOBANK COTE D''IVOIRE (ECOC)"
$data1[37,1] = 1
$data1[37,2] = 1
$data1[37,3] = -1.82
$data1[37,4] = -5.08
$data1[37,5] = "🟡 Observer"
$data1[37,6] = "👀 À surveiller"
$data1[38,0] = "SONATEL SN (SNTS)"
$data1[38,1] = 0
$data1[38,2] = 1
$data1[38,3] = -1.83
$data1[38,4] = -1.83
$data1[38,5] = "🟡 Observer"
$data1[38,6] = "➖ Neutre"
$data1[39,0] = "BANK OF AFRICA BN (BOAB)"
$data1[39,1] = 0
$data1[39,2] = 1
$data1[39,3] = -2
$data1[39,4] = -2
$data1[39,5] = "🟡 Observer"
$data1[39,6] = "➖ Neutre"
$data1[40,0] = "CFAO MOTORS CI (CFAC)"
$data1[40,1] = 0
$data1[40,2] = 1
$data1[40,3] = -2.21
$data1[40,4] = -2.21
$data1[40,5] = "🟡 Observer"
$data1[40,6] = "➖ Neutre"
$data1[41,0] = "BICI CI (BICC)"
$data1[41,1] = 0
$data1[41,2] = 1
$data1[41,3] = -2.6
$data1[41,4] = -2.6
$data1[41,5] = "🟡 Observer"
$data1[41,6] = "➖ Neutre"
$data1[42,0] = "BANK OF AFRICA BF (BOABF)"
$data1[42,1] = 0
$data1[42,2] = 1
$data1[42,3] = -2.86
$data1[42,4] = -2.86
$data1[42,5] = "🟡 Observer"
$data1[42,6] = "➖ Neutre"
$data1[43,0] = "SODE CI (SDCC)"
$data1[43,1] = 0
$data1[43,2] = 1
$data1[43,3] = -4.05
$data1[43,4] = -4.05
$data1[43,5] = "🟡 Observer"
$data1[43,6] = "➖ Neutre"
$data1[44,0] = "SOLIBRA CI (SLBC)"
$data1[44,1] = 1
$data1[44,2] = 2
$data1[44,3] = -4.87
$data1[44,4] = -4.87
$data1[44,5] = "🟡 Observer"
$data1[44,6] = "👀 À surveiller"
$data1[45,0] = "ECOBANK TRANS. INCORP. TG (ETIT)"
$data1[45,1] = 0
$data1[45,2] = 1
$data1[45,3] = -5.88
$data1[45,4] = -5.88
$data1[45,5] = "🟡 Observer"
$data1[45,6] = "➖ Neutre"
$data1[46,0] = "SETAO CI (STAC)"
$data1[46,1] = 0
$data1[46,2] = 1
$data1[46,3] = -6.09
$data1[46,4] = -6.09
$data1[46,5] = "🟡 Observer"
$data1[46,6] = "➖ Neutre"
$data1[47,0] = "SAPH CI (SPHC)"
$data1[47,1] = 0
$data1[47,2] = 1
$data1[47,3] = -7.46
$data1[47,4] = -7.46
$data1[47,5] = "🟡 Observer"
$data1[47,6] = "➖ Neutre"

$wsReco.Range("A2:G49").Value = $data1

# ---------------------------------------------------------------------------
# 2) "Top_YTD" sheet: rewrite the Progression YTD (%) values in rows 2-11.
# ---------------------------------------------------------------------------
$data2 = New-Object 'object[,]' 10,2
$data2[0,0] = "BRVM - SERVICES PUBLICS"
$data2[0,1] = 9325900.220000001
$data2[1,0] = "SUCRIVOIRE"
$data2[1,1] = 1360301.21
$data2[2,0] = "SAFCA CI"
$data2[2,1] = 394220
$data2[3,0] = "CFAO MOTORS CI"
$data2[3,1] = 360590.33
$data2[4,0] = "BRVM - AUTRES SECTEURS"
$data2[4,1] = 340873.21
$data2[5,0] = "NEI-CEDA CI"
$data2[5,1] = 231474
$data2[6,0] = "UNIWAX CI"
$data2[6,1] = 220037.08
$data2[7,0] = "SETAO CI"
$data2[7,1] = 196719.2
$data2[8,0] = "AIR LIQUIDE CI"
$data2[8,1] = 161200.16
$data2[9,0] = "BRVM - DISTRIBUTION"
$data2[9,1] = 49434.17

$wsYtd.Range("A2:B11").Value = $data2

Write-Output "BRVM data refreshed: Recommandations(A1:G49), Top_YTD(A1:B11)"
